# Ground_Data.xlsx - convert the CSV-derived sheet into the new layout:
#   - add "Id" / "Name" columns in front of the existing data
#   - rename the remaining headers (Start Depth/End Depth/Drop Items/
#     Sprite Addressable -> StartDepth/EndDepth/DropItems/SpriteAddressable)
#   - rebuild the table over the new A1:G5 range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing table first; we'll recreate it once the data is laid
# out in its final shape.
$lo = $ws.ListObjects.Item(1)
$lo.Delete()

# Write the data column-by-column (not row-by-row) so that shared-string
# insertion order matches how the sheet was actually authored.

# Column A: Id (numeric, no strings)
$ws.Cells.Item(1, 1).Value = "Id"
$ws.Cells.Item(2, 1).Value = 5001
$ws.Cells.Item(3, 1).Value = 5002
$ws.Cells.Item(4, 1).Value = 5003
$ws.Cells.Item(5, 1).Value = 5004

# Column B: Name
$ws.Cells.Item(1, 2).Value = "Name"
$ws.Cells.Item(2, 2).Value = "Ground_1"
$ws.Cells.Item(3, 2).Value = "Ground_2"
$ws.Cells.Item(4, 2).Value = "Ground_3"
$ws.Cells.Item(5, 2).Value = "Ground_4"

# Column C: StartDepth
$ws.Cells.Item(1, 3).Value = "StartDepth"
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(3, 3).Value = 6
$ws.Cells.Item(4, 3).Value = 11
$ws.Cells.Item(5, 3).Value = 16

# Column D: EndDepth
$ws.Cells.Item(1, 4).Value = "EndDepth"
$ws.Cells.Item(2, 4).Value = 5
$ws.Cells.Item(3, 4).Value = 10
$ws.Cells.Item(4, 4).Value = 15
$ws.Cells.Item(5, 4).Value = -1

# Column E: HP
$ws.Cells.Item(1, 5).Value = "HP"
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(3, 5).Value = 5
$ws.Cells.Item(4, 5).Value = 5
$ws.Cells.Item(5, 5).Value = 5

# Column F: DropItems
$ws.Cells.Item(1, 6).Value = "DropItems"
$ws.Cells.Item(2, 6).Value = "stone"
$ws.Cells.Item(3, 6).Value = "stone;iron"
$ws.Cells.Item(4, 6).Value = "stone;iron;gold"
$ws.Cells.Item(5, 6).Value = "diamond"

# Column G: SpriteAddressable
$ws.Cells.Item(1, 7).Value = "SpriteAddressable"
$ws.Cells.Item(2, 7).Value = "Light_Brown[Light_Brown]"
$ws.Cells.Item(3, 7).Value = "Deep_Brown[Deep_Brown]"
$ws.Cells.Item(4, 7).Value = "Lava_Earth[Lava_Earth]"
$ws.Cells.Item(5, 7).Value = "Light_Brown[Light_Brown]"

# Recreate the table over the full new range, keeping the original
# (Korean) table/display name.
$newlo = $ws.ListObjects.Add(1, $ws.Range("A1:G5"), $null, 1)
$newlo.Name = "표1"

# Match the saved view state: page setup + active selection.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$null = $ws.Range("G5").Select()
